# "Colocando header nos graficos"
# Adds a header label to column A (row 1) on each chart-data sheet,
# fixes accented Portuguese labels in column A, removes the now-unused
# header style from the data rows (keeping only the row-1 header style),
# drops the "Teto" row from the Emissoes sheet, and refreshes the
# "Custo Total" sheet headers/values.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheets 1-4 share the same row layout:
#   row1: (blank) | 2015 | 2030 | 2040 | 2050   -> add "Fonte/Tecnologia" to A1
#   rows 2-12: source/technology labels in column A
# ---------------------------------------------------------------
$sourceLabels = @{
    2  = "Hidro"
    3  = "Gás Natural"
    4  = "Carvão"
    5  = "Nuclear"
    6  = "Óleos Comb"
    7  = "Biomassa"
    8  = "Eólica"
    9  = "Solar"
    10 = "Outros"
    11 = "Pot. Compl."
    12 = "GD"
}

for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Give A1 the same header style already used by B1:E1, then set its text.
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    foreach ($row in $sourceLabels.Keys) {
        $cell = $ws.Cells.Item($row, 1)
        $cell.Value = $sourceLabels[$row]
        $cell.ClearFormats() | Out-Null
    }
}

# ---------------------------------------------------------------
# Sheet 5: "Emissoes Totais (MtCO2eq)"
#   row1: (blank) | 2015 | 2030 | 2040 | 2050   -> add "Período" to A1
#   row2: P Medio   -> P.Médio
#   row3: P Critico -> P.Crítico
#   row4: Teto      -> removed entirely
# ---------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("B1").Copy() | Out-Null
$ws5.Range("A1").PasteSpecial(-4122) | Out-Null
$ws5.Range("A1").Value = "Período"

$ws5.Cells.Item(2, 1).Value = "P.Médio"
$ws5.Cells.Item(2, 1).ClearFormats() | Out-Null

$ws5.Cells.Item(3, 1).Value = "P.Crítico"
$ws5.Cells.Item(3, 1).ClearFormats() | Out-Null

$ws5.Rows.Item(4).Delete() | Out-Null

# ---------------------------------------------------------------
# Sheet 6: "Custo Total (bilhões de R$)"
#   row1: Custo -> 2015, plus new A1 header "Tipo Expansão"
#   row2: Expansao Centralizada -> Expansão Centralizada, 722 -> 573
#   row3: Expansao por GD -> Expansão por GD, 65 -> 99
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("B1").Copy() | Out-Null
$ws6.Range("A1").PasteSpecial(-4122) | Out-Null
$ws6.Range("A1").Value = "Tipo Expansão"

# B1 must hold the literal text "2015" (not a number) while keeping its
# existing header style (s=1). Pasting the already-text "2015" value from
# another sheet's header (values-only paste) avoids Excel's automatic
# text->number coercion and keeps B1's current style untouched.
$wsFirst = $wb.Worksheets.Item(1)
$wsFirst.Range("B1").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4163) | Out-Null

$ws6.Cells.Item(2, 1).Value = "Expansão Centralizada"
$ws6.Cells.Item(2, 1).ClearFormats() | Out-Null
$ws6.Cells.Item(2, 2).Value = 573

$ws6.Cells.Item(3, 1).Value = "Expansão por GD"
$ws6.Cells.Item(3, 1).ClearFormats() | Out-Null
$ws6.Cells.Item(3, 2).Value = 99

Write-Host "Headers added and labels updated."
